# Saldo.xlsx update — refresh the "Export" balances extract with the new
# upload (accounts re-ranked by balance; a few accounts moved to brand new
# balances, a couple of stale small-balance rows are dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($rowNum, $conta, $nome, $saldo) {
    # Conta holds zero-padded account numbers ("004212581") — force text so
    # Excel doesn't coerce away the leading zeros. Flip the cell to text
    # format before the write (so the digit string is stored verbatim as
    # text, not re-parsed as a number), then clear the formatting override
    # back off so the cell is left with the sheet's default (unstyled)
    # look, matching the plain, un-styled data cells elsewhere in the
    # column.
    $contaCell = $ws.Cells.Item($rowNum, 1)
    $contaCell.NumberFormat = "@"
    $contaCell.Value = $conta
    $contaCell.ClearFormats()

    $ws.Cells.Item($rowNum, 2).Value = $nome
    $ws.Cells.Item($rowNum, 3).Value = $saldo
}

# --- Insertions (processed bottom-to-top so earlier row numbers used below
#     stay valid without re-deriving offsets) -------------------------------

# Row 5 (004207641 MAGALI 13693.23) is replaced by two new rows: overwrite
# it in place with the first (ALBERTO), then insert a second (MIRELLA)
# right after it.
Set-DataRow 5 "004376853" "ALBERTO" 10060.96
$ws.Rows(6).Insert()
Set-DataRow 6 "003553997" "MIRELLA" 6524.43

# Insert MAFALDA right before the existing row 4 (005366671 TATIANA).
$ws.Rows(4).Insert()
Set-DataRow 4 "004383190" "MAFALDA" 20570.49

# Insert JURACI right before the existing row 3 (004181486 ANDREA).
$ws.Rows(3).Insert()
Set-DataRow 3 "004342617" "JURACI" 37999.05

# Insert 7 new top-balance rows right before the existing row 2
# (004211922 CARLOS).
$ws.Rows("2:8").Insert()
Set-DataRow 2 "004212581" "MARIA"   129653.71
Set-DataRow 3 "004216401" "SUELY"   91133.06
Set-DataRow 4 "004387250" "MONICA"  70607.64
Set-DataRow 5 "005924958" "TIAGO"   50000
Set-DataRow 6 "004238436" "DIEGO"   44751.95
Set-DataRow 7 "004432935" "JOSE"    43772.06
Set-DataRow 8 "004381180" "HFR"     41410.9

# --- Deletions of now-superseded rows ---------------------------------
# The 10 rows inserted above push every following row down by 10, so the
# old small-balance duplicates (originally at row 80 and row 179) now sit
# at row 90 and row 189. Delete bottom-to-top.
$ws.Rows(189).Delete()  # old 004432935 JOSE 18.21
$ws.Rows(90).Delete()   # old 004376853 ALBERTO 60.96
